$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.795.54"
$ws.Range("E2").Value = "'  -0.90%  "
$ws.Range("D3").Value = "'2.089.53"
$ws.Range("E3").Value = "'  +1.89%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'246.30"
$ws.Range("E5").Value = "'  -0.90%  "
$ws.Range("D6").Value = "'0.651"
$ws.Range("E6").Value = "'  -2.07%  "
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("D8").Value = "'55.95"
$ws.Range("E8").Value = "'  -4.19%  "
$ws.Range("D9").Value = "'59.39"
$ws.Range("E9").Value = "'  -1.54%  "
$ws.Range("D10").Value = "'0.368"
$ws.Range("E10").Value = "'  -4.20%  "
$ws.Range("D11").Value = "'0.0769"
$ws.Range("E11").Value = "'  -1.88%  "
$ws.Range("E12").Value = "'  +1.34%  "
$ws.Range("D13").Value = "'14.88"
$ws.Range("E13").Value = "'  -5.32%  "
$ws.Range("D14").Value = "'0.881"
$ws.Range("E14").Value = "'  +6.00%  "
$ws.Range("D15").Value = "'2.386.70"
$ws.Range("E15").Value = "'  +1.47%  "
$ws.Range("E16").Value = "'  -3.78%  "
$ws.Range("D17").Value = "'2.096.55"
$ws.Range("E17").Value = "'  +2.16%  "
$ws.Range("D18").Value = "'36.789.51"
$ws.Range("E18").Value = "'  -0.93%  "
$ws.Range("D19").Value = "'17.50"
$ws.Range("E19").Value = "'  -1.85%  "
$ws.Range("D20").Value = "'73.12"
$ws.Range("E20").Value = "'  -2.28%  "
$ws.Range("D21").Value = "'0.0₃0879"
$ws.Range("D22").Value = "'5.48"
$ws.Range("E22").Value = "'  +1.88%  "
$ws.Range("E23").Value = "'  -0.45%  "
$ws.Range("E24").Value = "'  +0.03%  "
$ws.Range("D25").Value = "'2.42"
$ws.Range("E25").Value = "'  -2.35%  "
$ws.Range("D26").Value = "'9.97"
$ws.Range("E26").Value = "'  +6.73%  "
$ws.Range("D27").Value = "'2.17"
$ws.Range("E27").Value = "'  -0.49%  "
$ws.Range("D28").Value = "'168.27"
$ws.Range("E28").Value = "'  -0.49%  "
$ws.Range("D29").Value = "'21.06"
$ws.Range("E29").Value = "'  +4.93%  "
$ws.Range("E30").Value = "'  -0.68%  "
$ws.Range("D31").Value = "'5.33"
$ws.Range("E31").Value = "'  +10.49%  "
$ws.Range("D32").Value = "'1.20"
$ws.Range("E32").Value = "'  +7.15%  "
$ws.Range("D33").Value = "'4.73"
$ws.Range("E33").Value = "'  +3.95%  "
$ws.Range("E34").Value = "'  -1.61%  "
$ws.Range("D35").Value = "'2.39"
$ws.Range("E35").Value = "'  +5.58%  "
$ws.Range("E36").Value = "'  +0.09%  "
$ws.Range("D37").Value = "'1.83"
$ws.Range("E37").Value = "'  +4.42%  "
$ws.Range("D38").Value = "'0.0844"
$ws.Range("E38").Value = "'  -6.61%  "
$ws.Range("D39").Value = "'1.29"
$ws.Range("E39").Value = "'  -3.70%  "
$ws.Range("D40").Value = "'4.95"
$ws.Range("E40").Value = "'  -4.09%  "
$ws.Range("E41").Value = "'  +1.59%  "
$ws.Range("E42").Value = "'  -0.47%  "
$ws.Range("E44").Value = "'  -8.67%  "
$ws.Range("D45").Value = "'96.56"
$ws.Range("E45").Value = "'  +0.50%  "
$ws.Range("D46").Value = "'16.36"
$ws.Range("E46").Value = "'  -5.69%  "
$ws.Range("B47").Value = "'RenderToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'2.45"
$ws.Range("E47").Value = "'  -0.25%  "
$ws.Range("B48").Value = "'Maker"
$ws.Range("C48").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "'1.336.23"
$ws.Range("E48").Value = "'  +3.88%  "
$ws.Range("E49").Value = "'  +2.49%  "
$ws.Range("D50").Value = "'2.88"
$ws.Range("E50").Value = "'  -0.77%  "
$ws.Range("D51").Value = "'2.273.58"
$ws.Range("E51").Value = "'  +1.31%  "
